# Add two new years (2021年, 2022年) of data to the bottom of the table.
# Before the edit: data goes through row 10 (2020年). After: rows 11-12 added
# (2021年 fully populated, 2022年 only has the aggregate column AC filled in).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 11: 2021年 (fully populated)
# ---------------------------------------------------------------------------
$ws.Cells.Item(11, 1).Value = "2021年"
# Copy the label-cell formatting (bold, centered, bordered) from the row above
# so the new year label matches the rest of column A.
# -4122 == xlPasteFormats
$ws.Cells.Item(10, 1).Copy()
$ws.Cells.Item(11, 1).PasteSpecial(-4122)

$ws.Cells.Item(11, 2).Value = 1468.85
$ws.Cells.Item(11, 3).Value = 399.3
$ws.Cells.Item(11, 4).Value = 71.53
# Column E (其他采矿业) has no reported value for 2021年 - leave as blank marker.
$ws.Cells.Item(11, 5).Value = ""
$ws.Cells.Item(11, 5).Style = "Normal"
$ws.Cells.Item(11, 6).Value = 1345.04
$ws.Cells.Item(11, 7).Value = 2478.03
$ws.Cells.Item(11, 8).Value = 356.02
$ws.Cells.Item(11, 9).Value = 1385.88
$ws.Cells.Item(11, 10).Value = 301.63
$ws.Cells.Item(11, 11).Value = 330.89
$ws.Cells.Item(11, 12).Value = 249.62
$ws.Cells.Item(11, 13).Value = 4.44
$ws.Cells.Item(11, 14).Value = 531.64
$ws.Cells.Item(11, 15).Value = 1411.05
$ws.Cells.Item(11, 16).Value = 105.23
$ws.Cells.Item(11, 17).Value = 411.5
$ws.Cells.Item(11, 18).Value = 994.5599999999999
$ws.Cells.Item(11, 19).Value = 51.91
$ws.Cells.Item(11, 20).Value = 871.84
$ws.Cells.Item(11, 21).Value = 5.14
$ws.Cells.Item(11, 22).Value = 1273.23
$ws.Cells.Item(11, 23).Value = 99.3
$ws.Cells.Item(11, 24).Value = 400.95
$ws.Cells.Item(11, 25).Value = 2308.29
$ws.Cells.Item(11, 26).Value = 412.76
$ws.Cells.Item(11, 27).Value = 497
$ws.Cells.Item(11, 28).Value = 3.81
$ws.Cells.Item(11, 29).Value = 31774.13
$ws.Cells.Item(11, 30).Value = 927.4
$ws.Cells.Item(11, 31).Value = 539.66
$ws.Cells.Item(11, 32).Value = 2163.57
$ws.Cells.Item(11, 33).Value = 1454.47
$ws.Cells.Item(11, 34).Value = 362.68
$ws.Cells.Item(11, 35).Value = 597.8
$ws.Cells.Item(11, 36).Value = 24.58
$ws.Cells.Item(11, 37).Value = 1483.92
$ws.Cells.Item(11, 38).Value = 318.76
$ws.Cells.Item(11, 39).Value = 2941.54
$ws.Cells.Item(11, 40).Value = 252.6
$ws.Cells.Item(11, 41).Value = 680.29
$ws.Cells.Item(11, 42).Value = 1948.3
$ws.Cells.Item(11, 43).Value = 308.95

# ---------------------------------------------------------------------------
# Row 12: 2022年 (only the private-industry aggregate column AC is known so far)
# ---------------------------------------------------------------------------
$ws.Cells.Item(12, 1).Value = "2022年"
$ws.Cells.Item(10, 1).Copy()
$ws.Cells.Item(12, 1).PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(12, 29).Value = 26638

# All the other columns in row 12 are not yet reported - write them as
# present-but-blank cells (matching the rest of the still-empty row) rather
# than leaving the row sparse.
$blankCols = 2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,30,31,32,33,34,35,36,37,38,39,40,41,42,43
foreach ($col in $blankCols) {
    $ws.Cells.Item(12, $col).Value = ""
    $ws.Cells.Item(12, $col).Style = "Normal"
}
